$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the data analysis: B20 value was typo-prone when referenced by cell;
# replace the cell reference with the fixed constant 0.3042 in the formulas.
$ws.Range("C21").Formula = "=B21/0.3042*C20"
$ws.Range("C22:C29").Formula = "=B22/0.3042*C21"

# Daily entry: update the saved selection/active cell in the sheet view.
$ws.Activate()
$ws.Range("F23").Select()
